$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Append the new log row (row 6) for the repeat "Wanneer zijn jullie open?" test mail.
$ws.Cells.Item(6, 1).Value = "Wanneer zijn jullie open?"
$ws.Cells.Item(6, 2).Value = "mailmind.test@zohomail.eu"
$ws.Cells.Item(6, 3).Value = "Testmail #1: Wanneer zijn jullie open?"
$ws.Cells.Item(6, 4).Value = "Openingstijden / Locatie"
$ws.Cells.Item(6, 5).Value = "Beste klant,`nBedankt voor uw e-mail. Onze openingstijden zijn van maandag tot en met vrijdag van 9:00 tot 17:00 uur. We zijn gesloten in het weekend. Mocht u nog verdere vragen hebben, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Naam bedrijf]"
$ws.Cells.Item(6, 6).Value = "2025-06-27 22:24:11"
$ws.Cells.Item(6, 7).Value = "Ja"
$ws.Cells.Item(6, 8).Value = "Nee"
$ws.Cells.Item(6, 9).Value = "Ja"

# Extend the existing conditional-formatting blocks so they cover the new row too.
$ws.Range("D2:D5").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D6"))
$ws.Range("G2:G5").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G6"))
$ws.Range("H2:H5").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H6"))
$ws.Range("I2:I5").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I6"))

# Refresh the Dashboard summary count for "Openingstijden / Locatie" (2 -> 3).
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 3
